# Insert a new data row at row 189 (weekly Hortaliza/Fruta price record),
# which shifts the existing rows 189-249 down to 190-250.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("189:189").Insert()

# Populate the newly inserted row 189 with the new market record.
$ws.Range("A189").Value = 6
$ws.Range("B189").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C189").Value = "Metropolitana"
$ws.Range("D189").Value = 44809
$ws.Range("E189").Value = 13
$ws.Range("F189").Value = 100112022
$ws.Range("G189").Value = "Arveja Verde"
$ws.Range("H189").Value = "Perfection"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 300
$ws.Range("K189").Value = 35000
$ws.Range("L189").Value = 36000
$ws.Range("M189").Value = 35600
$ws.Range("N189").Value = "$/malla 25 kilos"
$ws.Range("O189").Value = "Provincia de Huasco"
$ws.Range("P189").Value = 1424
$ws.Range("Q189").Value = 25
$ws.Range("R189").Value = "Hortaliza"
